# Add a new, blank slide at the end of the deck (slide 4 / sldId 259).
#
# PowerPoint always writes a full slide skeleton (grpSpPr/xfrm, slide-level
# p14:creationId extLst, and clrMapOvr) even for an empty slide, so instead
# of using Slides.Add (which only yields a bare-bones <p:sld>), we duplicate
# the last existing slide -- which already uses the "Blank" layout shared by
# every slide in this deck -- and then strip out all of its shapes. That
# reproduces the exact skeleton PowerPoint emits while leaving the new slide
# completely empty, landing it at the end of the slide list in one step.

$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$source = $p.Slides.Item($lastIndex)

$newRange = $source.Duplicate()
$newSlide = $newRange.Item(1)

for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
    $newSlide.Shapes.Item($i).Delete()
}
